$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new numeric-looking price must stay text,
# matching how the rest of the "Price" column is stored as text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.403.00'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '3.495.37'
$ws.Range("E3").Value = '  -2.74%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '602.79'
$ws.Range("E5").Value = '  -3.31%  '
$ws.Range("D6").Value = '148.99'
$ws.Range("E6").Value = '  -4.70%  '
$ws.Range("D7").Value = '3.494.00'
$ws.Range("E7").Value = '  -2.69%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("D11").Value = '6.98'
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("E12").Value = '  -2.96%  '
$ws.Range("E13").Value = '  -3.75%  '
$ws.Range("D14").Value = '4.091.47'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").Value = '31.41'
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("D16").Value = '3.502.33'
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("D17").Value = '67.407.83'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '6.38'
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").Value = '15.07'
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("D21").Value = '445.12'
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("D22").Value = '9.07'
$ws.Range("E22").Value = '  -8.52%  '
$ws.Range("E23").Value = '  -3.47%  '
$ws.Range("D24").Value = '77.21'
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("D25").Value = '3.636.80'
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +8.00%  '
$ws.Range("D28").Value = '10.09'
$ws.Range("E28").Value = '  -6.04%  '
$ws.Range("D29").Value = '8.21'
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = '1.52'
$ws.Range("E32").Value = '  -6.74%  '
$ws.Range("E33").Value = '  +3.17%  '
$ws.Range("D34").Value = '25.56'
$ws.Range("E34").Value = '  -2.01%  '
$ws.Range("D35").Value = '3.480.36'
$ws.Range("E35").Value = '  -3.12%  '
$ws.Range("D36").Value = '6.04'
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("E37").Value = '  -4.79%  '
$ws.Range("D38").Value = '8.03'
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '178.15'
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").Value = '0.0877'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").Value = '5.37'
$ws.Range("E44").Value = '  -4.81%  '
$ws.Range("E45").Value = '  -3.21%  '
$ws.Range("D46").Value = '45.20'
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").Value = '27.60'
$ws.Range("E47").Value = '  -4.66%  '
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("D49").Value = '2.53'
$ws.Range("E49").Value = '  -2.33%  '
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("D51").Value = '0.987'
$ws.Range("E51").Value = '  -3.31%  '
